$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.529.54'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.012.05'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '508.45'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.50'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.59'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.530.85'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.48'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.96%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.572.56'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.24'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.016.31'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.94'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '329.46'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.70'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0916'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.75'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.33'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.33%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.61'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.11'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.73%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.41'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0675'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.046.87'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.58'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.96%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.651'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.224.77'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.12%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0239'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.40'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.85'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.64%  '
